$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (everything from old row 2 onward shifts down by one)
$ws.Rows.Item(2).Insert()

# The freshly inserted row inherits header-row formatting; reset it and copy the
# "data row" style (bold border + date number format) from the row right below,
# matching style s="2" used throughout column A.
$ws.Range("A2:E2").ClearFormats()
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2 (new first data point)
$ws.Cells.Item(2, 1).Value = 39400
$ws.Cells.Item(2, 2).Value = 2007
$ws.Cells.Item(2, 3).Value = 4.930115226412357
$ws.Cells.Item(2, 4).Value = 2008
$ws.Cells.Item(2, 5).ClearContents()

# Row 3
$ws.Cells.Item(3, 1).Value = 39765
$ws.Cells.Item(3, 2).Value = 2008
$ws.Cells.Item(3, 3).Value = 1.457587285166628
$ws.Cells.Item(3, 4).Value = 2009
$ws.Cells.Item(3, 5).ClearContents()

# Row 4
$ws.Cells.Item(4, 1).Value = 40130
$ws.Cells.Item(4, 2).Value = 2009
$ws.Cells.Item(4, 3).Value = -0.9140166223623458
$ws.Cells.Item(4, 4).Value = 2010
$ws.Cells.Item(4, 5).ClearContents()

# Row 5
$ws.Cells.Item(5, 1).Value = 40494
$ws.Cells.Item(5, 2).Value = 2010
$ws.Cells.Item(5, 3).Value = 2.585942866987878
$ws.Cells.Item(5, 4).Value = 2011
$ws.Cells.Item(5, 5).Value = 4.109775046142405

# Row 6
$ws.Cells.Item(6, 1).Value = 40862
$ws.Cells.Item(6, 2).Value = 2011
$ws.Cells.Item(6, 3).Value = 4.253963781362402
$ws.Cells.Item(6, 4).Value = 2012
$ws.Cells.Item(6, 5).Value = 2.863367440851095

# Row 7
$ws.Cells.Item(7, 1).Value = 41228
$ws.Cells.Item(7, 2).Value = 2012
$ws.Cells.Item(7, 3).Value = 1.752870900283909
$ws.Cells.Item(7, 4).Value = 2013
$ws.Cells.Item(7, 5).Value = 1.520397254708405

# Row 8
$ws.Cells.Item(8, 1).Value = 41592
$ws.Cells.Item(8, 2).Value = 2013
$ws.Cells.Item(8, 3).Value = -1.479696720105139
$ws.Cells.Item(8, 4).Value = 2014
$ws.Cells.Item(8, 5).Value = 2.503951807923066

# Row 9
$ws.Cells.Item(9, 1).Value = 41957
$ws.Cells.Item(9, 2).Value = 2014
$ws.Cells.Item(9, 3).Value = 3.900127535411246
$ws.Cells.Item(9, 4).Value = 2015
$ws.Cells.Item(9, 5).Value = 1.194160460927884

# Row 10
$ws.Cells.Item(10, 1).Value = 42321
$ws.Cells.Item(10, 2).Value = 2015
$ws.Cells.Item(10, 3).Value = 0.03947433952959933
$ws.Cells.Item(10, 4).Value = 2016
$ws.Cells.Item(10, 5).Value = 1.459149667419779

# Row 11
$ws.Cells.Item(11, 1).Value = 42689
$ws.Cells.Item(11, 2).Value = 2016
$ws.Cells.Item(11, 3).Value = 2.192778679161944
$ws.Cells.Item(11, 4).Value = 2017
$ws.Cells.Item(11, 5).Value = 1.586376095629216

# Row 12
$ws.Cells.Item(12, 1).Value = 43053
$ws.Cells.Item(12, 2).Value = 2017
$ws.Cells.Item(12, 3).Value = 3.40836448860673
$ws.Cells.Item(12, 4).Value = 2018
$ws.Cells.Item(12, 5).Value = 2.570658574505469

# Row 13
$ws.Cells.Item(13, 1).Value = 43418
$ws.Cells.Item(13, 2).Value = 2018
$ws.Cells.Item(13, 3).Value = 2.799070570134488
$ws.Cells.Item(13, 4).Value = 2019
$ws.Cells.Item(13, 5).Value = 2.479713128614147

# Row 14
$ws.Cells.Item(14, 1).Value = 43783
$ws.Cells.Item(14, 2).Value = 2019
$ws.Cells.Item(14, 3).Value = 4.195393191694419
$ws.Cells.Item(14, 4).Value = 2020
$ws.Cells.Item(14, 5).Value = 2.359935293525561

# Row 15
$ws.Cells.Item(15, 1).Value = 44159
$ws.Cells.Item(15, 2).Value = 2020
$ws.Cells.Item(15, 3).Value = 1.666553973046048
$ws.Cells.Item(15, 4).Value = 2021
$ws.Cells.Item(15, 5).Value = -0.4512719783814068

# Row 16
$ws.Cells.Item(16, 1).Value = 44525
$ws.Cells.Item(16, 2).Value = 2021
$ws.Cells.Item(16, 3).Value = 1.879266440112803
$ws.Cells.Item(16, 4).Value = 2022
$ws.Cells.Item(16, 5).Value = 1.081814991510499

# Row 17
$ws.Cells.Item(17, 1).Value = 44890
$ws.Cells.Item(17, 2).Value = 2022
$ws.Cells.Item(17, 3).Value = -2.620683231370946
$ws.Cells.Item(17, 4).Value = 2023
$ws.Cells.Item(17, 5).Value = -1.174318230871441

# Row 18
$ws.Cells.Item(18, 1).Value = 45254
$ws.Cells.Item(18, 2).Value = 2023
$ws.Cells.Item(18, 3).Value = -3.036556262700274
$ws.Cells.Item(18, 4).Value = 2024
$ws.Cells.Item(18, 5).Value = 0.07123445333143685

# Row 19
$ws.Cells.Item(19, 1).Value = 45618
$ws.Cells.Item(19, 2).Value = 2024
$ws.Cells.Item(19, 3).Value = -2.953443685011514
$ws.Cells.Item(19, 4).Value = 2025
$ws.Cells.Item(19, 5).Value = -1.196842846539037
